$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: date, variedad, prices, unit, price/kg, kg/unidad change. Origen (R9) stays "Región de O'Higgins".
$ws.Cells.Item(9,4).Value = 44902
$ws.Cells.Item(9,11).Value = "Castle Brite"
$ws.Cells.Item(9,14).Value = 15000
$ws.Cells.Item(9,15).Value = 16000
$ws.Cells.Item(9,16).Value = 15500
$ws.Cells.Item(9,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(9,19).Value = 1550
$ws.Cells.Item(9,20).Value = 10

# Row 10: date, variedad, prices, unit, price/kg, kg/unidad change. Origen (R10) stays "Región de O'Higgins".
$ws.Cells.Item(10,4).Value = 44902
$ws.Cells.Item(10,11).Value = "Castle Brite"
$ws.Cells.Item(10,14).Value = 13000
$ws.Cells.Item(10,15).Value = 13000
$ws.Cells.Item(10,16).Value = 13000
$ws.Cells.Item(10,17).Value = "$/caja 10 kilos"
$ws.Cells.Item(10,19).Value = 1300
$ws.Cells.Item(10,20).Value = 10

# Row 11: date, variedad, prices, unit, origen, price/kg, kg/unidad change.
$ws.Cells.Item(11,4).Value = 44559
$ws.Cells.Item(11,11).Value = "Modesto"
$ws.Cells.Item(11,14).Value = 19000
$ws.Cells.Item(11,15).Value = 20000
$ws.Cells.Item(11,16).Value = 19500
$ws.Cells.Item(11,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(11,18).Value = "Región de O'Higgins"
$ws.Cells.Item(11,19).Value = 1083
$ws.Cells.Item(11,20).Value = 18

# Row 12: date, variedad, calidad, prices, unit, origen, price/kg, kg/unidad change.
$ws.Cells.Item(12,4).Value = 44559
$ws.Cells.Item(12,11).Value = "Modesto"
$ws.Cells.Item(12,12).Value = "Segunda"
$ws.Cells.Item(12,14).Value = 18000
$ws.Cells.Item(12,15).Value = 18000
$ws.Cells.Item(12,16).Value = 18000
$ws.Cells.Item(12,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(12,18).Value = "Región de O'Higgins"
$ws.Cells.Item(12,19).Value = 1000
$ws.Cells.Item(12,20).Value = 18

# Row 13: new row (previously row 11 data, with its old date kept).
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13,3).Value = "Bíobío"
$ws.Cells.Item(13,4).Value = 44159
$ws.Cells.Item(13,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,5).Value = 8
$ws.Cells.Item(13,6).Value = "Fruta"
$ws.Cells.Item(13,7).Value = 100103
$ws.Cells.Item(13,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(13,9).Value = 100103003
$ws.Cells.Item(13,10).Value = "Damasco"
$ws.Cells.Item(13,11).Value = "Castle Brite"
$ws.Cells.Item(13,12).Value = "Primera"
$ws.Cells.Item(13,13).Value = 100
$ws.Cells.Item(13,14).Value = 14000
$ws.Cells.Item(13,15).Value = 15000
$ws.Cells.Item(13,16).Value = 14500
$ws.Cells.Item(13,17).Value = "$/caja 15 kilos"
$ws.Cells.Item(13,18).Value = "Región Metropolitana"
$ws.Cells.Item(13,19).Value = 967
$ws.Cells.Item(13,20).Value = 15

# Row 14: new row (previously row 12 data, with its old date kept).
$ws.Cells.Item(14,1).Value = 11
$ws.Cells.Item(14,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14,3).Value = "Bíobío"
$ws.Cells.Item(14,4).Value = 44875
$ws.Cells.Item(14,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14,5).Value = 8
$ws.Cells.Item(14,6).Value = "Fruta"
$ws.Cells.Item(14,7).Value = 100103
$ws.Cells.Item(14,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(14,9).Value = 100103003
$ws.Cells.Item(14,10).Value = "Damasco"
$ws.Cells.Item(14,11).Value = "Castle Brite"
$ws.Cells.Item(14,12).Value = "Primera"
$ws.Cells.Item(14,13).Value = 50
$ws.Cells.Item(14,14).Value = 31000
$ws.Cells.Item(14,15).Value = 32000
$ws.Cells.Item(14,16).Value = 31400
$ws.Cells.Item(14,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(14,18).Value = "Provincia de Limarí"
$ws.Cells.Item(14,19).Value = 3140
$ws.Cells.Item(14,20).Value = 10
